$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row changes ---
# A new "半成品等级" (semi-finished product level) column is inserted at C,
# pushing the old "半成品类型" (C) into D, and the old separate
# "枪械工专属"/"护甲工专属" profession columns (D/E) are consolidated into a
# single "专属职业" (exclusive profession) column at E. "获取方式" (F) and
# "配方" (G) keep their text.
$ws.Range("C1").Value = "半成品等级"
$ws.Range("D1").Value = "半成品类型"
$ws.Range("E1").Value = "专属职业"
$ws.Range("G1").Value = "配方"

# --- Column width ---
# Column C now shares the (wider) formatting previously used only by column B.
$ws.Columns("C:C").ColumnWidth = $ws.Columns("B:B").ColumnWidth

# --- Selection ---
# The author ended up with column C selected.
$ws.Columns("C:C").Select()
